# Insert a new data row at row 42 (pushes existing rows 42-121 down to 43-122)
# and populate it with the new record from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(42).Insert()

$ws.Range("A42").Value = 6
$ws.Range("B42").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44725
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100104
$ws.Range("H42").Value = "Frutos de pepita"
$ws.Range("I42").Value = 100104003
$ws.Range("J42").Value = "Membrillo"
$ws.Range("K42").Value = "Champion"
$ws.Range("L42").Value = "Especial"
$ws.Range("M42").Value = 8
$ws.Range("N42").Value = 300000
$ws.Range("O42").Value = 300000
$ws.Range("P42").Value = 300000
$ws.Range("Q42").Value = "$/bins (450 kilos)"
$ws.Range("R42").Value = "Región de O'Higgins"
$ws.Range("S42").Value = 667
$ws.Range("T42").Value = 450
